# IRC_Chat_Program.docx edit script
# Rewrites the TODO bullet list, the STRUCTURE bullet list (bookmark move +
# new Debug.h bullet), merges a few split runs in the DOCUMENTATION section,
# and appends the new m_ConsoleLog / OutLogDump() documentation paragraphs.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ParaXml($inner) {
    return "<w:p $wNs>$inner</w:p>"
}

# ---------------------------------------------------------------------
# TODO: list (numId=2) - paragraphs 4..13 (COM 1-based paragraph index)
# ---------------------------------------------------------------------
$pPrTodo = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'

# 1) "DEBUG logger..." -> "MAIN will house a TRY function for IRC_PROGRAM..."
$inner1 = $pPrTodo + '<w:r><w:rPr><w:b/></w:rPr><w:t>MAIN</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> will house a TRY function for </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>IRC_PROGRAM</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> that starts the whole program, and the CATCH will be a crasher for log reports.</w:t></w:r>'
$d.Paragraphs(4).Range.InsertXML((New-ParaXml $inner1))

# 2) "MAIN will house..." -> "ArgV[1] will be reserved for debugging."
$inner2 = $pPrTodo + '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>ArgV</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>[</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>1] will be reserved for debugging.</w:t></w:r>'
$d.Paragraphs(5).Range.InsertXML((New-ParaXml $inner2))

# 3) "Add debugger outputs for DEBUG" -> "OPTIONALLY: SSL/SSH? Do not save passwords as plaintext."
$inner3 = $pPrTodo + '<w:r><w:t>OPTIONALLY: SSL/SSH?</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Do not save passwords as plaintext.</w:t></w:r>'
$d.Paragraphs(6).Range.InsertXML((New-ParaXml $inner3))

# 4) "If DEBUG is set to 0..." -> "User/Password storage in XML format."
$inner4 = $pPrTodo + '<w:r><w:t>User/Password storage in XML format.</w:t></w:r>'
$d.Paragraphs(7).Range.InsertXML((New-ParaXml $inner4))

# 5) "ArgV[1]..." -> "SERVER program created."
$inner5 = $pPrTodo + '<w:r><w:rPr><w:b/></w:rPr><w:t>SERVER</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> program</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>created.</w:t></w:r>'
$d.Paragraphs(8).Range.InsertXML((New-ParaXml $inner5))

# 6) "OPTIONALLY: SSL/SSH?..." -> "CLIENT program created."
$inner6 = $pPrTodo + '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">CLIENT </w:t></w:r>' +
    '<w:r><w:t>program created.</w:t></w:r>'
$d.Paragraphs(9).Range.InsertXML((New-ParaXml $inner6))

# 7) "User/Password storage..." -> "Documentation of the program please."
$inner7 = $pPrTodo + '<w:r><w:t>Documentation of the program please.</w:t></w:r>'
$d.Paragraphs(10).Range.InsertXML((New-ParaXml $inner7))

# 8) "SERVER program created." -> "Proper error crashing. Use OutLogDump()"
$inner8 = $pPrTodo + '<w:r><w:t xml:space="preserve">Proper error crashing. Use </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>OutLogDump</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>)</w:t></w:r>'
$d.Paragraphs(11).Range.InsertXML((New-ParaXml $inner8))

# 9) "CLIENT program created." -> "Set up network factory to pass off server connections" + bookmark
$inner9 = $pPrTodo + '<w:r><w:t>Set up network factory to pass off server connections</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$d.Paragraphs(12).Range.InsertXML((New-ParaXml $inner9))

# 10) "Documentation of the program please." paragraph is removed entirely
#     (the old bullet list had 10 items, the new one has 9).
$d.Paragraphs(13).Range.Delete()

Write-Output "TODO list done"

# ---------------------------------------------------------------------
# STRUCTURE: list (numId=4)
# ---------------------------------------------------------------------
$pPrStructure = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>'

# Find the "CLIENT will request and authenticate..." bullet (the bookmark
# used to sit mid-run splitting "chi" + "ld thread."); merge back into one
# run and drop the bookmark (it now lives on the new TODO bullet instead).
$clientPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith('CLIENT will request and authenticate')) {
        $clientPara = $i
        break
    }
}
$innerClient = $pPrStructure + '<w:r><w:rPr><w:b/></w:rPr><w:t>CLIENT</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> will request and authenticate. A separate thread will process SEND/RECEIVE and leave main thread solely for user input and doing SEND requests to child thread.</w:t></w:r>'
$d.Paragraphs($clientPara).Range.InsertXML((New-ParaXml $innerClient))

# Insert the new "Debug.h houses all debug functions..." bullet right after
# the last numId=4 bullet ("CLIENT will be able to submit chat commands...").
$lastStructPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith('CLIENT will be able to submit chat commands')) {
        $lastStructPara = $i
        break
    }
}
$innerDebugH = $pPrStructure + '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Debug.h</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>houses all debug functions and error log dumping. Please use this throughout the program.</w:t></w:r>'
$insertPos = $d.Paragraphs($lastStructPara).Range.End - 1
$d.Range($insertPos, $insertPos).InsertXML((New-ParaXml $innerDebugH))

Write-Output "STRUCTURE list done"

# ---------------------------------------------------------------------
# DOCUMENTATION: section - merge split "WriteImportantMessage(" / "const " /
# "string message)" runs (and the Informational twin) into single runs.
# ---------------------------------------------------------------------
function Find-ParaStartingWith($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return $null
}

$pPr720 = '<w:pPr><w:ind w:left="720"/></w:pPr>'

$impPara = Find-ParaStartingWith "`tDebug::WriteImportantMessage"
$innerImp = $pPr720 + '<w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Debug::</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>WriteImportantMessage(const string message)</w:t></w:r>'
$d.Paragraphs($impPara).Range.InsertXML((New-ParaXml $innerImp))

$infoPara = Find-ParaStartingWith "`tDebug::WriteInformationalMessage"
$innerInfo = $pPr720 + '<w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Debug::</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>WriteInformationalMessage(const string message)</w:t></w:r>'
$d.Paragraphs($infoPara).Range.InsertXML((New-ParaXml $innerInfo))

Write-Output "Message signature merges done"

# ---------------------------------------------------------------------
# Append the m_ConsoleLog / OutLogDump() documentation right after the
# "This sets the level..." paragraph, reusing the blank <w:ind left=360/>
# paragraph that already sits there for the first new line.
# ---------------------------------------------------------------------
$pPr360 = '<w:pPr><w:ind w:left="360"/></w:pPr>'

$blankPara = Find-ParaStartingWith "This sets the level for the program"
$blankPara = $blankPara + 1   # the empty <w:ind left=360/> paragraph right after it

$innerConsoleLog = $pPr360 + '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>Debug::</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>m_ConsoleLog</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$d.Paragraphs($blankPara).Range.InsertXML((New-ParaXml $innerConsoleLog))

$innerVectorDesc = $pPr360 + '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t>Vector pointer to all messages outputted by the debugger.</w:t></w:r>'
$innerOutLogDump = $pPr360 + '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>Debug::</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>OutLogDump</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$innerOutLogDumpDesc = $pPr360 + '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">Creates a text file IRC_Program_Log.txt with </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>m_ConsoleLog</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>'

$newParasXml = (New-ParaXml $innerVectorDesc) + (New-ParaXml $innerOutLogDump) + (New-ParaXml $innerOutLogDumpDesc)
$insertPos = $d.Paragraphs($blankPara).Range.End - 1
$d.Range($insertPos, $insertPos).InsertXML($newParasXml)

Write-Output "ConsoleLog/OutLogDump documentation added"
